# "report preparation and give-data link for employee"
#
# Updates the "issues" sheet (Sheet2):
#  - row 7  changes from "unique employee id check" to "link for employee"
#  - a new row 8 is inserted with three notes about a "give data" report/link
#  - "unique employee id check" moves down to row 9
#  - a new row 10 holds a note about the "give data" action
#  - the existing backlog rows (data reset.../back buttons/.../what should
#    happen...) shift down to rows 12-16, leaving row 11 blank
#  - two new notes are placed further down the sheet, in J20 and K21
#  - column A/B get wider, and the active selection moves to A11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("issues")

# Make room: push the existing "data reset.../.../what should happen..."
# block (old rows 9-13) down by 3 rows so it starts at row 12, opening up
# rows 9, 10 and 11 right after the current row 7/8.
$ws.Rows("9:11").Insert()

# New strings must be introduced in the same order they first appear in the
# finished sheet, so the shared-string table is appended to in that order:
# A7, J20, K21, A8, B8, C8, (A9 reuses the existing "unique employee id
# check" string), A10.
$ws.Range("A7").Value = "link for employee"
$ws.Range("J20").Value = "Techincal"
$ws.Range("K21").Value = "change required defaulters name as it is also used for report"
$ws.Range("A8").Value = "check over all functionality"
$ws.Range("B8").Value = "defaulters email ids"
$ws.Range("C8").Value = "link copy for action, to share with employee"
$ws.Range("A9").Value = "unique employee id check"
$ws.Range("A10").Value = "give data : action is displayed for all pax"

# Widen the first two columns to fit the new, longer text.
$ws.Columns("A").ColumnWidth = 39.166666666666664
$ws.Columns("B").ColumnWidth = 16.053385416666668

# Match the saved selection/active cell.
$ws.Range("A11").Select()
